# Refresh the demo/test claim data used by the automation scripts.
# (Mirrors the "Fixed demo scripts issue" commit: new ClaimsID / dates /
#  tracking numbers / waybill numbers / pickup-dropoff labels.)

$wb = $excel.ActiveWorkbook

# --- ShipmentInformation sheet ------------------------------------------
$wsShipment = $wb.Worksheets.Item("ShipmentInformation")
$wsShipment.Range("C2").Value = "PickUp25"
$wsShipment.Range("K2").Value = "DropOff715"

# --- Input sheet ---------------------------------------------------------
$wsInput = $wb.Worksheets.Item("Input")
$wsInput.Range("T2").Value = ""

# Force text storage ("@") before assigning numeric/date-looking strings so
# they stay text (t="s") instead of being auto-converted to a date/number.
$wsInput.Range("B3").NumberFormat = "@"
$wsInput.Range("B3").Value = "08-24-2021"

$wsInput.Range("T3").NumberFormat = "@"
$wsInput.Range("T3").Value = "57535724"

$wsInput.Range("W3").Value = "FCT879658849361985536"
$wsInput.Range("X3").Value = "FCTEST1003540"

# --- ClaimDetail sheet -----------------------------------------------------
$wsClaim = $wb.Worksheets.Item("ClaimDetail")
$wsClaim.Range("A3").NumberFormat = "@"
$wsClaim.Range("A3").Value = "57535724"

$wsClaim.Range("B3").NumberFormat = "@"
$wsClaim.Range("B3").Value = "08-24-2021"
